$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Formed the consolidated report: fill in the "Absent" column (H) for the
# remaining days so the attendance totals are complete.
$ws.Range("H11").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
